$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = '025da00d-c7d8-4f5b-9d8b-af717f607614_fila_9.png'
$ws.Range("B5").Value = 'Roman Alfonso Grios Boza'

$ws.Range("A6").Value = '1d7d3ca2-ae91-448b-a30f-af7eab2a3978_fila_15.png'
$ws.Range("B6").Value = 'Orlando Mauricio Guevara'

$ws.Range("A7").Value = '1ea081eb-fb66-47f3-89e9-c8d075c539d8_fila_7.png'
$ws.Range("B7").Value = 'Cristina Jozabed Carvajal'

$ws.Range("A8").Value = '3088def7-058f-48af-80af-065076c246c4_fila_2.png'
$ws.Range("B8").Value = 'Isabella Dompe Estrada'

$ws.Range("A9").Value = '39a55e95-dcfd-45f2-aabc-e9b0290f16dd_fila_10.png'
$ws.Range("B9").Value = 'Abraham Silva Ampre'

$ws.Range("A10").Value = '3b6de5df-f0f9-4db5-8b6f-0e7e1a36e18d_fila_4.png'
$ws.Range("B10").Value = 'Yadder Fernando Torres'

$ws.Range("A11").Value = '4342cc05-bbbd-4d09-a2a2-06a6308c1337_fila_13.png'
$ws.Range("B11").Value = 'Marlon Josue Gonzales Cano'

$ws.Range("A12").Value = '4ea9709c-9a55-4a66-8562-d88a0007d197_fila_6.png'
$ws.Range("B12").Value = 'Ronier Jose Rivera'

$ws.Range("A13").Value = '6580d43c-5690-42f2-80f6-4d378a94affe_fila_14.png'
$ws.Range("B13").Value = 'Angel Isaac Alvarez Quiñonez'

$ws.Range("A14").Value = '7dd4e43e-d4d2-4fee-977e-b13b6f1c0891_fila_11.png'
$ws.Range("B14").Value = 'Eduardo Domingo Zeledon Merca'

$ws.Range("A15").Value = '7e15a2d1-7b5f-454e-9c6f-d7df58f942c8_fila_3.png'
$ws.Range("B15").Value = 'Bryan Alexander Cano'

$ws.Range("A16").Value = '8facc7e6-bab2-46a0-a3d5-c6673c096089_fila_8.png'
$ws.Range("B16").Value = 'David Orlando Mena Valverd'

$ws.Range("A17").Value = 'c4c25ca5-7f6c-4c1a-a605-cc3c09be6af7_fila_5.png'
$ws.Range("B17").Value = 'Erick Espinoza'

$ws.Range("A18").Value = 'ca020dc0-339a-4fee-b903-5359cade4990_fila_1.png'
$ws.Range("B18").Value = 'Hotep Antonio Ruiz Lezama'

$ws.Range("A19").Value = 'cb9afc98-07fa-433d-b7e5-573c05fb955c_fila_12.png'
$ws.Range("B19").Value = 'José Danilo Suárez'
